$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Formula = '="43556201    "'
$ws.Cells.Item(1,2).Value = 'CHECOLI GONZALEZ PAULA JESUS  '
$ws.Cells.Item(1,3).Value = 'OBRA SOCIAL DE LA CONFEDERACION DE OBREROS Y EMPLEADOS MUNICIPALES ARGENTINA'

$ws.Cells.Item(2,1).Formula = '="22887465    "'
$ws.Cells.Item(2,2).Value = 'LINDNER MARTA ROSANA          '
$ws.Cells.Item(2,3).Value = 'OBRA SOCIAL DEL PERSONAL DE LA INDUSTRIA DEL CAUCHO ; INSTITUTO NACIONAL DE SERVICIOS SOCIALES PARA JUBILADOS Y PENSIONADOS'

$ws.Cells.Item(3,1).Formula = '="21706616    "'
$ws.Cells.Item(3,2).Value = 'THEILL JOSE MARIA             '
$ws.Cells.Item(3,3).Value = 'O.S.P. BUENOS AIRES (IOMA)'

$ws.Cells.Item(4,1).Formula = '="17273992    "'
$ws.Cells.Item(4,2).Value = 'NICASTRO SILVIA BEATRIZ       '
$ws.Cells.Item(4,3).Value = 'Programa Federal Incluir Salud (PFIS) '

$ws.Cells.Item(5,1).Formula = '="45453562    "'
$ws.Cells.Item(5,2).Value = 'BARRAZA AGUSTIN ALBERTO       '
$ws.Cells.Item(5,3).Value = 'OBRA SOCIAL DEL PERSONAL RURAL Y ESTIBADORES DE LA REPUBLICA ARGENTINA'

$ws.Cells.Item(6,1).Formula = '="92083965    "'
$ws.Cells.Item(6,2).Value = 'MARTINEZ VILLALBA MARIA ANGELA'
$ws.Cells.Item(6,3).Value = 'INSTITUTO NACIONAL DE SERVICIOS SOCIALES PARA JUBILADOS Y PENSIONADOS'

$ws.Cells.Item(7,1).Formula = '="22654192    "'
$ws.Cells.Item(7,2).Value = 'JAEGER LORENA ELISABETH       '
$ws.Cells.Item(7,3).Value = 'OBRA SOCIAL DE DOCENTES PARTICULARES ; O.S.P. BUENOS AIRES (IOMA)'

$ws.Cells.Item(8,1).Formula = '="54594940    "'
$ws.Cells.Item(8,2).Value = 'DIAZORIANA MARTINA            '
$ws.Cells.Item(8,3).Value = 'OBRA SOCIAL UNION PERSONAL DE LA UNION DEL  PERSONAL CIVIL DE LA NACION'

$ws.Cells.Item(9,1).Formula = '="20235877    "'
$ws.Cells.Item(9,2).Value = 'WARNHOLTZ JORGE FEDERICO      '
$ws.Cells.Item(9,3).Value = 'OBRA SOCIAL DE CAPITANES DE ULTRAMAR Y OFICIALES DE LA MARINA MERCANTE'

$ws.Cells.Item(10,1).Formula = '="16427478    "'
$ws.Cells.Item(10,2).Value = 'GHIORZI PABLO HUGO            '
$ws.Cells.Item(10,3).Value = 'O.S.P. BUENOS AIRES (IOMA)'

$ws.Cells.Item(11,1).Formula = '="57273676    "'
$ws.Cells.Item(11,2).Value = 'ROIBON NICOLAS                '
$ws.Cells.Item(11,3).Value = 'O.S.P. CORRIENTES (IOSCOR) ; OBRA SOCIAL DE EJECUTIVOS Y DEL PERSONAL DE DIRECCION DE EMPRESAS'

$ws.Cells.Item(12,1).Formula = '="36206894    "'
$ws.Cells.Item(12,2).Value = 'CASCO LAUREANO JOSE           '
$ws.Cells.Item(12,3).Value = 'OBRA SOCIAL DEL PERSONAL DE LA INDUSTRIA MADERERA'

$ws.Cells.Item(13,1).Formula = '="37182967    "'
$ws.Cells.Item(13,2).Value = 'PIZZORNO MARIA BELEN          '
$ws.Cells.Item(13,3).Value = 'INSTITUTO NACIONAL DE SERVICIOS SOCIALES PARA JUBILADOS Y PENSIONADOS'

$ws.Cells.Item(14,1).Formula = '="29508358    "'
$ws.Cells.Item(14,2).Value = 'GOMEZ HUGO ROBERTO            '
$ws.Cells.Item(14,3).Value = 'OBRA SOCIAL DE LA UNION OBRERA METALURGICA DE LA REPUBLICA ARGENTINA'

$ws.Cells.Item(15,1).Formula = '="31556222    "'
$ws.Cells.Item(15,2).Value = 'FERNANDEZ MARCOS MATIAS       '
$ws.Cells.Item(15,3).Value = 'OBRA SOCIAL DEL PERSONAL DE BARRACAS DE LANAS, CUEROS Y ANEXOS'

$ws.Cells.Item(16,1).Formula = '="20184875    "'
$ws.Cells.Item(16,2).Value = 'PAZ MARCOS FLAVIO             '
$ws.Cells.Item(16,3).Value = 'O.S.P. BUENOS AIRES (IOMA)'

$ws.Cells.Item(17,1).Formula = '="39847551    "'
$ws.Cells.Item(17,2).Value = 'AUTERI NAARA MILENA           '
$ws.Cells.Item(17,3).Value = 'OBRA SOCIAL DEL PERSONAL LADRILLERO'

$ws.Range("A1:A17").Copy()
$ws.Range("A1:A17").PasteSpecial(-4163)
$excel.CutCopyMode = 0